$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: CV_Summary  (row order + values change)
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("CV_Summary")

$wsSummary.Range("A2").Value = "10hz"
$wsSummary.Range("B2").Value = 0.83125
$wsSummary.Range("C2").Value = 0.04238956239453293
$wsSummary.Range("D2").Value = 0.90625
$wsSummary.Range("E2").Value = 0.78125
$wsSummary.Range("F2").Value = 19.60977554251545
$wsSummary.Range("G2").Value = 160
$wsSummary.Range("H2").Value = 89

$wsSummary.Range("A3").Value = "30hz"
$wsSummary.Range("B3").Value = 0.825
$wsSummary.Range("C3").Value = 0.025
$wsSummary.Range("D3").Value = 0.875
$wsSummary.Range("E3").Value = 0.8125
$wsSummary.Range("F3").Value = 32.99998680000527
$wsSummary.Range("G3").Value = 160
$wsSummary.Range("H3").Value = 90

$wsSummary.Range("A4").Value = "40hz"
$wsSummary.Range("B4").Value = 0.7875
$wsSummary.Range("C4").Value = 0.0125
$wsSummary.Range("D4").Value = 0.8125
$wsSummary.Range("E4").Value = 0.78125
$wsSummary.Range("F4").Value = 62.99994960004032
$wsSummary.Range("G4").Value = 160
$wsSummary.Range("H4").Value = 89

$wsSummary.Range("A5").Value = "20hz"
$wsSummary.Range("B5").Value = 0.75
$wsSummary.Range("C5").Value = 0.0625
$wsSummary.Range("D5").Value = 0.875
$wsSummary.Range("E5").Value = 0.71875
$wsSummary.Range("F5").Value = 11.99999808000031
$wsSummary.Range("G5").Value = 160
$wsSummary.Range("H5").Value = 89

# ---------------------------------------------------------------
# Sheet: CV_Scores_Detail  (only column C accuracy values change)
# ---------------------------------------------------------------
$wsDetail = $wb.Worksheets.Item("CV_Scores_Detail")

$wsDetail.Range("C2").Value = 0.90625
$wsDetail.Range("C3").Value = 0.78125
$wsDetail.Range("C4").Value = 0.84375
$wsDetail.Range("C5").Value = 0.8125
$wsDetail.Range("C6").Value = 0.8125
$wsDetail.Range("C7").Value = 0.71875
$wsDetail.Range("C8").Value = 0.71875
$wsDetail.Range("C9").Value = 0.875
$wsDetail.Range("C10").Value = 0.71875
$wsDetail.Range("C11").Value = 0.71875
$wsDetail.Range("C12").Value = 0.8125
$wsDetail.Range("C13").Value = 0.8125
$wsDetail.Range("C14").Value = 0.8125
$wsDetail.Range("C15").Value = 0.8125
$wsDetail.Range("C16").Value = 0.875
$wsDetail.Range("C17").Value = 0.78125
$wsDetail.Range("C18").Value = 0.78125
$wsDetail.Range("C19").Value = 0.8125
$wsDetail.Range("C20").Value = 0.78125
$wsDetail.Range("C21").Value = 0.78125

# ---------------------------------------------------------------
# Sheet: Analysis_Info
# ---------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Analysis_Info")

$wsInfo.Range("B3").Value = "10hz"
$wsInfo.Range("B4").Value = "20hz"
$wsInfo.Range("B5").Value = "40hz"
$wsInfo.Range("B6").Value = "2025-10-02 20:31:00"
